$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44267
$ws.Cells.Item(2, 12).Value = 'Primera'
$ws.Cells.Item(2, 13).Value = 50
$ws.Cells.Item(2, 14).Value = 15000
$ws.Cells.Item(2, 15).Value = 16000
$ws.Cells.Item(2, 16).Value = 15400
$ws.Cells.Item(2, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(2, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(2, 19).Value = 962
$ws.Cells.Item(2, 20).Value = 16
$ws.Cells.Item(3, 4).Value = 44249
$ws.Cells.Item(3, 12).Value = 'Primera'
$ws.Cells.Item(3, 13).Value = 110
$ws.Cells.Item(3, 14).Value = 16000
$ws.Cells.Item(3, 15).Value = 16000
$ws.Cells.Item(3, 16).Value = 16000
$ws.Cells.Item(3, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(3, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(3, 19).Value = 1000
$ws.Cells.Item(3, 20).Value = 16
$ws.Cells.Item(4, 4).Value = 44266
$ws.Cells.Item(4, 12).Value = 'Especial'
$ws.Cells.Item(4, 13).Value = 50
$ws.Cells.Item(4, 14).Value = 18000
$ws.Cells.Item(4, 15).Value = 19000
$ws.Cells.Item(4, 16).Value = 18400
$ws.Cells.Item(4, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(4, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(4, 19).Value = 1150
$ws.Cells.Item(4, 20).Value = 16
$ws.Cells.Item(5, 4).Value = 44266
$ws.Cells.Item(5, 12).Value = 'Primera'
$ws.Cells.Item(5, 13).Value = 120
$ws.Cells.Item(5, 14).Value = 15000
$ws.Cells.Item(5, 15).Value = 15000
$ws.Cells.Item(5, 16).Value = 15000
$ws.Cells.Item(5, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(5, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(5, 19).Value = 938
$ws.Cells.Item(5, 20).Value = 16
$ws.Cells.Item(6, 4).Value = 44266
$ws.Cells.Item(6, 12).Value = 'Segunda'
$ws.Cells.Item(6, 13).Value = 20
$ws.Cells.Item(6, 14).Value = 10000
$ws.Cells.Item(6, 15).Value = 10000
$ws.Cells.Item(6, 16).Value = 10000
$ws.Cells.Item(6, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(6, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(6, 19).Value = 625
$ws.Cells.Item(6, 20).Value = 16
$ws.Cells.Item(7, 4).Value = 44221
$ws.Cells.Item(7, 12).Value = 'Primera'
$ws.Cells.Item(7, 13).Value = 30
$ws.Cells.Item(7, 14).Value = 25000
$ws.Cells.Item(7, 15).Value = 25000
$ws.Cells.Item(7, 16).Value = 25000
$ws.Cells.Item(7, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(7, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(7, 19).Value = 1389
$ws.Cells.Item(7, 20).Value = 18
$ws.Cells.Item(8, 4).Value = 44251
$ws.Cells.Item(8, 12).Value = 'Primera'
$ws.Cells.Item(8, 13).Value = 55
$ws.Cells.Item(8, 14).Value = 16000
$ws.Cells.Item(8, 15).Value = 16000
$ws.Cells.Item(8, 16).Value = 16000
$ws.Cells.Item(8, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(8, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(8, 19).Value = 1000
$ws.Cells.Item(8, 20).Value = 16
$ws.Cells.Item(9, 4).Value = 44277
$ws.Cells.Item(9, 12).Value = 'Primera'
$ws.Cells.Item(9, 13).Value = 90
$ws.Cells.Item(9, 14).Value = 16000
$ws.Cells.Item(9, 15).Value = 16000
$ws.Cells.Item(9, 16).Value = 16000
$ws.Cells.Item(9, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(9, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(9, 19).Value = 1000
$ws.Cells.Item(9, 20).Value = 16
$ws.Cells.Item(10, 4).Value = 44237
$ws.Cells.Item(10, 12).Value = 'Primera'
$ws.Cells.Item(10, 13).Value = 20
$ws.Cells.Item(10, 14).Value = 25000
$ws.Cells.Item(10, 15).Value = 25000
$ws.Cells.Item(10, 16).Value = 25000
$ws.Cells.Item(10, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(10, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(10, 19).Value = 1562
$ws.Cells.Item(10, 20).Value = 16
$ws.Cells.Item(11, 4).Value = 44253
$ws.Cells.Item(11, 12).Value = 'Primera'
$ws.Cells.Item(11, 13).Value = 55
$ws.Cells.Item(11, 14).Value = 16000
$ws.Cells.Item(11, 15).Value = 16000
$ws.Cells.Item(11, 16).Value = 16000
$ws.Cells.Item(11, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(11, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(11, 19).Value = 1000
$ws.Cells.Item(11, 20).Value = 16
$ws.Cells.Item(12, 4).Value = 44323
$ws.Cells.Item(12, 12).Value = 'Segunda'
$ws.Cells.Item(12, 13).Value = 30
$ws.Cells.Item(12, 14).Value = 14000
$ws.Cells.Item(12, 15).Value = 14000
$ws.Cells.Item(12, 16).Value = 14000
$ws.Cells.Item(12, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(12, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(12, 19).Value = 875
$ws.Cells.Item(12, 20).Value = 16
$ws.Cells.Item(13, 4).Value = 44238
$ws.Cells.Item(13, 12).Value = 'Primera'
$ws.Cells.Item(13, 13).Value = 65
$ws.Cells.Item(13, 14).Value = 14000
$ws.Cells.Item(13, 15).Value = 14000
$ws.Cells.Item(13, 16).Value = 14000
$ws.Cells.Item(13, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(13, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(13, 19).Value = 875
$ws.Cells.Item(13, 20).Value = 16
$ws.Cells.Item(14, 4).Value = 44327
$ws.Cells.Item(14, 12).Value = 'Primera'
$ws.Cells.Item(14, 13).Value = 35
$ws.Cells.Item(14, 14).Value = 15000
$ws.Cells.Item(14, 15).Value = 15000
$ws.Cells.Item(14, 16).Value = 15000
$ws.Cells.Item(14, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(14, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(14, 19).Value = 938
$ws.Cells.Item(14, 20).Value = 16
$ws.Cells.Item(15, 4).Value = 44418
$ws.Cells.Item(15, 12).Value = 'Primera'
$ws.Cells.Item(15, 13).Value = 40
$ws.Cells.Item(15, 14).Value = 30000
$ws.Cells.Item(15, 15).Value = 30000
$ws.Cells.Item(15, 16).Value = 30000
$ws.Cells.Item(15, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(15, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(15, 19).Value = 1875
$ws.Cells.Item(15, 20).Value = 16
$ws.Cells.Item(16, 4).Value = 44284
$ws.Cells.Item(16, 12).Value = 'Primera'
$ws.Cells.Item(16, 13).Value = 25
$ws.Cells.Item(16, 14).Value = 18000
$ws.Cells.Item(16, 15).Value = 18000
$ws.Cells.Item(16, 16).Value = 18000
$ws.Cells.Item(16, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(16, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(16, 19).Value = 1125
$ws.Cells.Item(16, 20).Value = 16
$ws.Cells.Item(17, 4).Value = 44244
$ws.Cells.Item(17, 12).Value = 'Primera'
$ws.Cells.Item(17, 13).Value = 25
$ws.Cells.Item(17, 14).Value = 14000
$ws.Cells.Item(17, 15).Value = 14000
$ws.Cells.Item(17, 16).Value = 14000
$ws.Cells.Item(17, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(17, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(17, 19).Value = 875
$ws.Cells.Item(17, 20).Value = 16
$ws.Cells.Item(18, 4).Value = 44330
$ws.Cells.Item(18, 12).Value = 'Especial'
$ws.Cells.Item(18, 13).Value = 35
$ws.Cells.Item(18, 14).Value = 30000
$ws.Cells.Item(18, 15).Value = 30000
$ws.Cells.Item(18, 16).Value = 30000
$ws.Cells.Item(18, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(18, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(18, 19).Value = 1500
$ws.Cells.Item(18, 20).Value = 20
$ws.Cells.Item(19, 4).Value = 44265
$ws.Cells.Item(19, 12).Value = 'Primera'
$ws.Cells.Item(19, 13).Value = 50
$ws.Cells.Item(19, 14).Value = 16000
$ws.Cells.Item(19, 15).Value = 16000
$ws.Cells.Item(19, 16).Value = 16000
$ws.Cells.Item(19, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(19, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(19, 19).Value = 1000
$ws.Cells.Item(19, 20).Value = 16
$ws.Cells.Item(20, 4).Value = 44315
$ws.Cells.Item(20, 12).Value = 'Primera'
$ws.Cells.Item(20, 13).Value = 105
$ws.Cells.Item(20, 14).Value = 15000
$ws.Cells.Item(20, 15).Value = 16000
$ws.Cells.Item(20, 16).Value = 15619
$ws.Cells.Item(20, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(20, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(20, 19).Value = 976
$ws.Cells.Item(20, 20).Value = 16
$ws.Cells.Item(21, 4).Value = 44294
$ws.Cells.Item(21, 12).Value = 'Primera'
$ws.Cells.Item(21, 13).Value = 80
$ws.Cells.Item(21, 14).Value = 14000
$ws.Cells.Item(21, 15).Value = 16000
$ws.Cells.Item(21, 16).Value = 15000
$ws.Cells.Item(21, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(21, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(21, 19).Value = 938
$ws.Cells.Item(21, 20).Value = 16
$ws.Cells.Item(22, 4).Value = 44427
$ws.Cells.Item(22, 12).Value = 'Primera'
$ws.Cells.Item(22, 13).Value = 65
$ws.Cells.Item(22, 14).Value = 40000
$ws.Cells.Item(22, 15).Value = 40000
$ws.Cells.Item(22, 16).Value = 40000
$ws.Cells.Item(22, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(22, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(22, 19).Value = 2500
$ws.Cells.Item(22, 20).Value = 16
$ws.Cells.Item(23, 4).Value = 44319
$ws.Cells.Item(23, 12).Value = 'Primera'
$ws.Cells.Item(23, 13).Value = 50
$ws.Cells.Item(23, 14).Value = 15000
$ws.Cells.Item(23, 15).Value = 15000
$ws.Cells.Item(23, 16).Value = 15000
$ws.Cells.Item(23, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(23, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(23, 19).Value = 938
$ws.Cells.Item(23, 20).Value = 16
$ws.Cells.Item(24, 4).Value = 44273
$ws.Cells.Item(24, 12).Value = 'Primera'
$ws.Cells.Item(24, 13).Value = 55
$ws.Cells.Item(24, 14).Value = 14000
$ws.Cells.Item(24, 15).Value = 14000
$ws.Cells.Item(24, 16).Value = 14000
$ws.Cells.Item(24, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(24, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(24, 19).Value = 875
$ws.Cells.Item(24, 20).Value = 16
$ws.Cells.Item(25, 4).Value = 44236
$ws.Cells.Item(25, 12).Value = 'Primera'
$ws.Cells.Item(25, 13).Value = 45
$ws.Cells.Item(25, 14).Value = 25000
$ws.Cells.Item(25, 15).Value = 25000
$ws.Cells.Item(25, 16).Value = 25000
$ws.Cells.Item(25, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(25, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(25, 19).Value = 1562
$ws.Cells.Item(25, 20).Value = 16
$ws.Cells.Item(26, 4).Value = 44235
$ws.Cells.Item(26, 12).Value = 'Primera'
$ws.Cells.Item(26, 13).Value = 55
$ws.Cells.Item(26, 14).Value = 25000
$ws.Cells.Item(26, 15).Value = 25000
$ws.Cells.Item(26, 16).Value = 25000
$ws.Cells.Item(26, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(26, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(26, 19).Value = 1562
$ws.Cells.Item(26, 20).Value = 16
$ws.Cells.Item(27, 4).Value = 44280
$ws.Cells.Item(27, 12).Value = 'Especial'
$ws.Cells.Item(27, 13).Value = 80
$ws.Cells.Item(27, 14).Value = 20000
$ws.Cells.Item(27, 15).Value = 20000
$ws.Cells.Item(27, 16).Value = 20000
$ws.Cells.Item(27, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(27, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(27, 19).Value = 1250
$ws.Cells.Item(27, 20).Value = 16
$ws.Cells.Item(28, 4).Value = 44274
$ws.Cells.Item(28, 12).Value = 'Especial'
$ws.Cells.Item(28, 13).Value = 35
$ws.Cells.Item(28, 14).Value = 16000
$ws.Cells.Item(28, 15).Value = 16000
$ws.Cells.Item(28, 16).Value = 16000
$ws.Cells.Item(28, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(28, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(28, 19).Value = 1000
$ws.Cells.Item(28, 20).Value = 16
$ws.Cells.Item(29, 4).Value = 44274
$ws.Cells.Item(29, 12).Value = 'Primera'
$ws.Cells.Item(29, 13).Value = 40
$ws.Cells.Item(29, 14).Value = 14000
$ws.Cells.Item(29, 15).Value = 14000
$ws.Cells.Item(29, 16).Value = 14000
$ws.Cells.Item(29, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(29, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(29, 19).Value = 875
$ws.Cells.Item(29, 20).Value = 16
$ws.Cells.Item(30, 4).Value = 44320
$ws.Cells.Item(30, 12).Value = 'Primera'
$ws.Cells.Item(30, 13).Value = 50
$ws.Cells.Item(30, 14).Value = 14000
$ws.Cells.Item(30, 15).Value = 14000
$ws.Cells.Item(30, 16).Value = 14000
$ws.Cells.Item(30, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(30, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(30, 19).Value = 875
$ws.Cells.Item(30, 20).Value = 16
$ws.Cells.Item(31, 4).Value = 44258
$ws.Cells.Item(31, 12).Value = 'Primera'
$ws.Cells.Item(31, 13).Value = 65
$ws.Cells.Item(31, 14).Value = 16000
$ws.Cells.Item(31, 15).Value = 16000
$ws.Cells.Item(31, 16).Value = 16000
$ws.Cells.Item(31, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(31, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(31, 19).Value = 889
$ws.Cells.Item(31, 20).Value = 18
$ws.Cells.Item(32, 4).Value = 44435
$ws.Cells.Item(32, 12).Value = 'Especial'
$ws.Cells.Item(32, 13).Value = 30
$ws.Cells.Item(32, 14).Value = 45000
$ws.Cells.Item(32, 15).Value = 45000
$ws.Cells.Item(32, 16).Value = 45000
$ws.Cells.Item(32, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(32, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(32, 19).Value = 2500
$ws.Cells.Item(32, 20).Value = 18
$ws.Cells.Item(33, 4).Value = 44306
$ws.Cells.Item(33, 12).Value = 'Especial'
$ws.Cells.Item(33, 13).Value = 50
$ws.Cells.Item(33, 14).Value = 22000
$ws.Cells.Item(33, 15).Value = 22000
$ws.Cells.Item(33, 16).Value = 22000
$ws.Cells.Item(33, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(33, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(33, 19).Value = 1375
$ws.Cells.Item(33, 20).Value = 16
$ws.Cells.Item(34, 4).Value = 44301
$ws.Cells.Item(34, 12).Value = 'Especial'
$ws.Cells.Item(34, 13).Value = 55
$ws.Cells.Item(34, 14).Value = 22000
$ws.Cells.Item(34, 15).Value = 22000
$ws.Cells.Item(34, 16).Value = 22000
$ws.Cells.Item(34, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(34, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(34, 19).Value = 1375
$ws.Cells.Item(34, 20).Value = 16
$ws.Cells.Item(35, 4).Value = 44301
$ws.Cells.Item(35, 12).Value = 'Primera'
$ws.Cells.Item(35, 13).Value = 85
$ws.Cells.Item(35, 14).Value = 19000
$ws.Cells.Item(35, 15).Value = 19000
$ws.Cells.Item(35, 16).Value = 19000
$ws.Cells.Item(35, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(35, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(35, 19).Value = 1188
$ws.Cells.Item(35, 20).Value = 16
$ws.Cells.Item(36, 4).Value = 44321
$ws.Cells.Item(36, 12).Value = 'Primera'
$ws.Cells.Item(36, 13).Value = 100
$ws.Cells.Item(36, 14).Value = 15000
$ws.Cells.Item(36, 15).Value = 15000
$ws.Cells.Item(36, 16).Value = 15000
$ws.Cells.Item(36, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(36, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(36, 19).Value = 938
$ws.Cells.Item(36, 20).Value = 16
$ws.Cells.Item(37, 4).Value = 44425
$ws.Cells.Item(37, 12).Value = 'Primera'
$ws.Cells.Item(37, 13).Value = 35
$ws.Cells.Item(37, 14).Value = 40000
$ws.Cells.Item(37, 15).Value = 40000
$ws.Cells.Item(37, 16).Value = 40000
$ws.Cells.Item(37, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(37, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(37, 19).Value = 2500
$ws.Cells.Item(37, 20).Value = 16
$ws.Cells.Item(38, 4).Value = 44270
$ws.Cells.Item(38, 12).Value = 'Primera'
$ws.Cells.Item(38, 13).Value = 55
$ws.Cells.Item(38, 14).Value = 10000
$ws.Cells.Item(38, 15).Value = 10000
$ws.Cells.Item(38, 16).Value = 10000
$ws.Cells.Item(38, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(38, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(38, 19).Value = 625
$ws.Cells.Item(38, 20).Value = 16
$ws.Cells.Item(39, 4).Value = 44389
$ws.Cells.Item(39, 12).Value = 'Especial'
$ws.Cells.Item(39, 13).Value = 30
$ws.Cells.Item(39, 14).Value = 30000
$ws.Cells.Item(39, 15).Value = 30000
$ws.Cells.Item(39, 16).Value = 30000
$ws.Cells.Item(39, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(39, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(39, 19).Value = 1875
$ws.Cells.Item(39, 20).Value = 16
$ws.Cells.Item(40, 4).Value = 44326
$ws.Cells.Item(40, 12).Value = 'Primera'
$ws.Cells.Item(40, 13).Value = 100
$ws.Cells.Item(40, 14).Value = 15000
$ws.Cells.Item(40, 15).Value = 15000
$ws.Cells.Item(40, 16).Value = 15000
$ws.Cells.Item(40, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(40, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(40, 19).Value = 938
$ws.Cells.Item(40, 20).Value = 16
$ws.Cells.Item(41, 4).Value = 44432
$ws.Cells.Item(41, 12).Value = 'Especial'
$ws.Cells.Item(41, 13).Value = 30
$ws.Cells.Item(41, 14).Value = 45000
$ws.Cells.Item(41, 15).Value = 45000
$ws.Cells.Item(41, 16).Value = 45000
$ws.Cells.Item(41, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(41, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(41, 19).Value = 2500
$ws.Cells.Item(41, 20).Value = 18
